$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.475.15'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '1.904.09'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4793'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4067'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08075'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.002'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.40%  '
$ws.Range("D12").Value = '1.894.98'
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.955'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.074'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06684'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001033'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").Value = '29.486.52'
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.542'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.59%  '
$ws.Range("E24").Value = '  -2.08%  '
$ws.Range("D25").Value = '2.121.72'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.089'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.097'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.49'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.040'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09510'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.444'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.394'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.542'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.98%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02251'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.77%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06075'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5881'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("E40").Value = '  -5.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1846'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.420'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.278'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07797'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5530'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.924'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '113.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.2944'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.02%  '
